# Restored from revision: update the "From" value of rule R30 (row 10)
# in the Sample Project rules table from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
